$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $style = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $style
}

$ws.Range("D2").Value = '52.169.79'
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").Value = '2.977.48'
$ws.Range("E3").Value = '  +1.39%  '
$ws.Range("E4").Value = '  +0.04%  '
Set-TextValue $ws.Range("D5") '354.37'
$ws.Range("E5").Value = '  +0.56%  '
Set-TextValue $ws.Range("D6") '107.28'
$ws.Range("E6").Value = '  -4.53%  '
Set-TextValue $ws.Range("D7") '0.563'
$ws.Range("E7").Value = '  +0.62%  '
$ws.Range("E8").Value = '  +0.00%  '
Set-TextValue $ws.Range("D9") '0.615'
$ws.Range("E9").Value = '  -2.04%  '
Set-TextValue $ws.Range("D10") '38.27'
$ws.Range("E10").Value = '  -2.94%  '
$ws.Range("E11").Value = '  +1.35%  '
$ws.Range("E12").Value = '  -3.65%  '
Set-TextValue $ws.Range("D13") '19.23'
$ws.Range("E13").Value = '  -3.93%  '
$ws.Range("D14").Value = '3.443.53'
$ws.Range("E14").Value = '  +1.38%  '
Set-TextValue $ws.Range("D15") '7.63'
$ws.Range("E15").Value = '  -2.62%  '
$ws.Range("D16").Value = '2.976.03'
$ws.Range("E16").Value = '  +1.70%  '
$ws.Range("E17").Value = '  +0.99%  '
$ws.Range("D18").Value = '52.195.84'
$ws.Range("E18").Value = '  +0.40%  '
Set-TextValue $ws.Range("D19") '3.48'
$ws.Range("E19").Value = '  +4.96%  '
Set-TextValue $ws.Range("D20") '7.49'
$ws.Range("E20").Value = '  -2.07%  '
Set-TextValue $ws.Range("D21") '13.61'
$ws.Range("E21").Value = '  -4.65%  '
$ws.Range("E22").Value = '  -1.29%  '
Set-TextValue $ws.Range("D23") '69.60'
Set-TextValue $ws.Range("D24") '263.96'
$ws.Range("E24").Value = '  -1.99%  '
$ws.Range("E25").Value = '  -1.43%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("B27").Value = 'Filecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D27") '7.64'
$ws.Range("E27").Value = '  +2.93%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range("D28") '26.80'
$ws.Range("E28").Value = '  -0.67%  '
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("E30").Value = '  +0.58%  '
$ws.Range("E31").Value = '  -3.00%  '
Set-TextValue $ws.Range("D32") '6.13'
$ws.Range("E32").Value = '  -1.34%  '
Set-TextValue $ws.Range("D33") '36.37'
$ws.Range("E33").Value = '  -3.03%  '
$ws.Range("E34").Value = '  -4.38%  '
Set-TextValue $ws.Range("D35") '50.77'
$ws.Range("E35").Value = '  -4.03%  '
$ws.Range("E36").Value = '  -2.27%  '
$ws.Range("E37").Value = '  +0.01%  '
Set-TextValue $ws.Range("D38") '3.20'
$ws.Range("E38").Value = '  -3.05%  '
Set-TextValue $ws.Range("D39") '17.93'
$ws.Range("E39").Value = '  -5.16%  '
Set-TextValue $ws.Range("D40") '1.97'
$ws.Range("E40").Value = '  -3.95%  '
Set-TextValue $ws.Range("D41") '2.72'
$ws.Range("E41").Value = '  +0.29%  '
$ws.Range("E42").Value = '  -0.42%  '
Set-TextValue $ws.Range("D43") '22.78'
$ws.Range("E43").Value = '  -1.71%  '
Set-TextValue $ws.Range("D44") '121.49'
$ws.Range("E44").Value = '  +8.25%  '
$ws.Range("E45").Value = '  -3.20%  '
$ws.Range("D46").Value = '2.117.47'
$ws.Range("E46").Value = '  -2.51%  '
$ws.Range("E47").Value = '  -4.31%  '
$ws.Range("E48").Value = '  -7.16%  '
Set-TextValue $ws.Range("D49") '0.240'
$ws.Range("E49").Value = '  -3.60%  '
Set-TextValue $ws.Range("D50") '0.0335'
$ws.Range("E50").Value = '  -2.82%  '
Set-TextValue $ws.Range("D51") '0.938'
$ws.Range("E51").Value = '  -0.10%  '
